$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.4471118299071225
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.5128108231650346
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.5142868200499384
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.3515300111981882
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.6667864272184959
$ws.Range("D7").Value = 0.5086072680161079
$ws.Range("D8").Value = 0.5248751753221683
$ws.Range("D9").Value = 0.5444834636029661
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0.4923343050435468
$ws.Range("D11").Value = 0.4537622034155825
$ws.Range("D12").Value = 0.5045749746481826
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0.4652021096909096
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0.5410010079298252
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.3555360560952532
$ws.Range("D16").Value = 0.5071481808419751
$ws.Range("D17").Value = 0.5122179597512977
$ws.Range("D18").Value = 0.4580393340350952
$ws.Range("D19").Value = 0.4650436118196076
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0.3366350046586943
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0.4902810119079558
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0.6086035789939667
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.7161874614834476
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0.511246202900055
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0.3692236982531322
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.4877370591280287
$ws.Range("D27").Value = 0.4970606232485276
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0.3912964603780546
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0.5096097894449901
$ws.Range("D30").Value = 0.5758828619899141
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 0.526239801506054
$ws.Range("D32").Value = 0.8293458909726055
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0.4630516867642296
$ws.Range("D34").Value = 0.4703648402918553
$ws.Range("D35").Value = 0.517538917014704
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 0.5383897165248632
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0.4984627778484001
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 0.4572142749627597
$ws.Range("D39").Value = 0.4428643489351282
$ws.Range("D40").Value = 0.4760624041373397
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 0.5023452308338077
$ws.Range("D42").Value = 0.5092943616825444
$ws.Range("D43").Value = 0.5234000653433046
$ws.Range("D44").Value = 0.5075424536570992
$ws.Range("D45").Value = 0.5075089696045735
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0.435398248250259
$ws.Range("D47").Value = 0.5012338855090729
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 0.6103260771724672
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 0.4471817559693818
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 0.4489785187772473
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 0.5094724402526877
$ws.Range("D52").Value = 0.28378730957202
$ws.Range("D53").Value = 0.5000821257894242
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 0.497518593356151
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 0.5374257846788126
$ws.Range("D56").Value = 0.4996938952755273
$ws.Range("D57").Value = 0.4500169551299417
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 0.4920383581699381
$ws.Range("D59").Value = 0.5091866002208123
$ws.Range("D60").Value = 0.4943503796385175
$ws.Range("D61").Value = 0.52190942274153
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = 0.5227540861106106
$ws.Range("D63").Value = 0.4729942202324009
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 0.5271650490963673
$ws.Range("D65").Value = 0.4958423544867442
$ws.Range("D66").Value = 0.512144208391502
$ws.Range("D67").Value = 0.3889865748370847
$ws.Range("D68").Value = 0.4973815685686857
$ws.Range("D69").Value = 0.5205294481876545
$ws.Range("D70").Value = 0.5442024818414894
$ws.Range("D71").Value = 0.4361085757248254
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 0.5106021877771861
$ws.Range("D73").Value = 0.498210102699182
$ws.Range("D74").Value = 0.4836213295320703
$ws.Range("C75").Value = 1
$ws.Range("D75").Value = 0.5141895309331896
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = 0.5239902968770835
$ws.Range("D77").Value = 0.4954826510045377
$ws.Range("D78").Value = 0.5076496448636435
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 0.5097464919733984
